# Update "想去人数" (interested count) figures across the three sheets
# that contain data rows: 展览 (sheet1), 演出 (sheet2), 全部类型 (sheet4).
# 本地生活 (sheet3) only has a header row and needs no changes.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 179
$ws1.Range("F5").Value  = 5065
$ws1.Range("F8").Value  = 16
$ws1.Range("F13").Value = 1414
$ws1.Range("F14").Value = 3741
$ws1.Range("F15").Value = 416
$ws1.Range("F19").Value = 2733
$ws1.Range("F20").Value = 137
$ws1.Range("F21").Value = 32
$ws1.Range("F25").Value = 69
$ws1.Range("F26").Value = 9
$ws1.Range("F29").Value = 278
$ws1.Range("F30").Value = 47

# --- 演出 sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 46

# --- 全部类型 sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 179
$ws4.Range("F5").Value  = 46
$ws4.Range("F6").Value  = 5065
$ws4.Range("F9").Value  = 16
$ws4.Range("F14").Value = 1414
$ws4.Range("F15").Value = 3741
$ws4.Range("F16").Value = 416
$ws4.Range("F20").Value = 2733
$ws4.Range("F21").Value = 137
$ws4.Range("F22").Value = 32
$ws4.Range("F26").Value = 69
$ws4.Range("F27").Value = 9
$ws4.Range("F30").Value = 278
$ws4.Range("F31").Value = 47
